# Nexial unitTest_vision.xlsx update
# - base: remove "clear(variables)" from the #system!F (base) list
# - external: append "terminate(programName)" to the #system!J (external) list
# - io: insert "assertPath(path)" into the #system!L (io) list (alphabetical position)
# - web: rename "assertAttributeContains(...)" -> "assertAttributeContain(...)" and
#        insert "saveSelectedText(var,locator)" / "saveSelectedValue(var,locator)"
#        into the #system!Z (web) list (alphabetical position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # '#system' sheet

function Find-Index($list, $value) {
    for ($i = 0; $i -lt $list.Count; $i++) {
        if ($list[$i] -eq $value) {
            return $i
        }
    }
    return -1
}

function Read-Column($sheet, $colIdx, $startRow, $endRow) {
    $result = New-Object System.Collections.ArrayList
    for ($r = $startRow; $r -le $endRow; $r++) {
        $result.Add($sheet.Cells.Item($r, $colIdx).Value2) | Out-Null
    }
    return $result
}

function Write-Column($sheet, $colIdx, $startRow, $list) {
    $r = $startRow
    foreach ($v in $list) {
        $sheet.Cells.Item($r, $colIdx).Value2 = $v
        $r = $r + 1
    }
    return $r
}

# ---------------------------------------------------------------------------
# Column F = "base" (currently F2:F40, 39 items)
# ---------------------------------------------------------------------------
$colF = 6
$oldLastF = 40
$listF = Read-Column $ws $colF 2 $oldLastF
$listF.Remove("clear(variables)") | Out-Null
$nextRowF = Write-Column $ws $colF 2 $listF
for ($r = $nextRowF; $r -le $oldLastF; $r++) {
    $ws.Cells.Item($r, $colF).Value2 = $null
}
$newLastF = $nextRowF - 1

# ---------------------------------------------------------------------------
# Column J = "external" (currently J2:J5, 4 items) -> append new entry
# ---------------------------------------------------------------------------
$colJ = 10
$oldLastJ = 5
$newLastJ = $oldLastJ + 1
$ws.Cells.Item($newLastJ, $colJ).Value2 = "terminate(programName)"

# ---------------------------------------------------------------------------
# Column L = "io" (currently L2:L29, 28 items) -> insert new entry
# ---------------------------------------------------------------------------
$colL = 12
$oldLastL = 29
$listL = Read-Column $ws $colL 2 $oldLastL
$insertIdxL = Find-Index $listL "assertReadableFile(file,minByte)"
$listL.Insert($insertIdxL, "assertPath(path)") | Out-Null
$nextRowL = Write-Column $ws $colL 2 $listL
$newLastL = $nextRowL - 1

# ---------------------------------------------------------------------------
# Column Z = "web" (currently Z2:Z135, 134 items)
#   -> rename "assertAttributeContains(...)" to "assertAttributeContain(...)"
#   -> insert 2 new entries before "saveTableAsCsv(...)"
# ---------------------------------------------------------------------------
$colZ = 26
$oldLastZ = 135
$listZ = Read-Column $ws $colZ 2 $oldLastZ

$renameIdxZ = Find-Index $listZ "assertAttributeContains(locator,attrName,contains)"
$listZ[$renameIdxZ] = "assertAttributeContain(locator,attrName,contains)"

$insertIdxZ = Find-Index $listZ "saveTableAsCsv(locator,nextPageLocator,file)"
$listZ.Insert($insertIdxZ, "saveSelectedValue(var,locator)") | Out-Null
$listZ.Insert($insertIdxZ, "saveSelectedText(var,locator)") | Out-Null

$nextRowZ = Write-Column $ws $colZ 2 $listZ
$newLastZ = $nextRowZ - 1

# ---------------------------------------------------------------------------
# Update defined names to reflect the new ranges
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$$newLastF"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$$newLastJ"
$wb.Names.Item("io").RefersTo = "='#system'!`$L`$2:`$L`$$newLastL"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$$newLastZ"

Write-Host "base ->" $wb.Names.Item("base").RefersTo
Write-Host "external ->" $wb.Names.Item("external").RefersTo
Write-Host "io ->" $wb.Names.Item("io").RefersTo
Write-Host "web ->" $wb.Names.Item("web").RefersTo
